# change predict income algorithm
# Roll the yearly income-statement columns (D:H) one fiscal period forward:
#   D <- old E, E <- old F, F <- old G, G <- old H, H <- new period.
# This drops the oldest fiscal year (1396/12) and its publish-date column,
# and appends the newest fiscal year (1401/12) with its publish date(s).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header rows: fiscal period (row 8) and publish date (row 9) ----
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

$ws.Range("D9").Value = "1399-04-11 (7)"
$ws.Range("E9").Value = "1400-04-09 (8)"
$ws.Range("F9").Value = "1401-04-01 (8)"
$ws.Range("G9").Value = "1402-02-30 (8)"
$ws.Range("H9").Value = "1402-02-30 (2)"

# ---- Data rows: shift D:H left by one period, fill H with the new value ----
function Shift-Row($row, $newH) {
    $e = $ws.Range("E$row").Value2
    $f = $ws.Range("F$row").Value2
    $g = $ws.Range("G$row").Value2
    $h = $ws.Range("H$row").Value2
    $ws.Range("D$row").Value = $e
    $ws.Range("E$row").Value = $f
    $ws.Range("F$row").Value = $g
    $ws.Range("G$row").Value = $h
    $ws.Range("H$row").Value = $newH
}

Shift-Row 11 36311
Shift-Row 12 -24190
Shift-Row 13 12121
Shift-Row 14 -2530
Shift-Row 15 "-"
Shift-Row 16 -24
Shift-Row 17 9566
Shift-Row 18 -306
Shift-Row 19 138
Shift-Row 20 9398
Shift-Row 21 -1130
Shift-Row 22 8268
Shift-Row 23 "-"
Shift-Row 24 8268
Shift-Row 25 0
Shift-Row 26 12706
Shift-Row 27 0
